$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("screenTitles")
$ws4.Range("A10").Value = "Shop"
$ws4.Range("B10").Value = "shop"
$ws4.Range("A1:B1").Select()

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws6.Name = "Films"

$ws6.Range("A1").Value = "objectID"
$ws6.Range("B1").Value = "name_nl"
$hdr = $ws6.Range("A1:B1")
$hdr.Font.Bold = $true
$hdr.Interior.ThemeColor = 0
$hdr.Interior.TintAndShade = -0.249977111117893
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

$ws6.Range("A2").Value = "RentedMovieCategorie"
$ws6.Range("A3").Value = "RentedMovieName"

$ws6.Range("B2").Value = "AKGDCC1 NL Branch"
$ws6.Range("B3").Value = "AKG DCC VOD1"

$ws6.Columns.Item(1).ColumnWidth = 22.28515625
$ws6.Columns.Item(2).ColumnWidth = 20

$ws6.Range("A2").Select()
$ws6.Activate()
